$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels (B1:E1 renamed to reflect TRUE/FALSE naming)
$ws.Range("B1").Value = "FALSE_count"
$ws.Range("C1").Value = "FALSE_percent"
$ws.Range("D1").Value = "TRUE_count"
$ws.Range("E1").Value = "TRUE_percent"

# Row 2: 5-10km
$ws.Range("A2").Value = "5-10km"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 5.555555555555555
$ws.Range("D2").Value = 17
$ws.Range("E2").Value = 22.97297297297298
$ws.Range("F2").Value = 9
$ws.Range("G2").Value = 23.07692307692308

# Row 3: <5km
$ws.Range("A3").Value = "<5km"
$ws.Range("B3").Value = 12
$ws.Range("C3").Value = 66.66666666666666
$ws.Range("D3").Value = 33
$ws.Range("E3").Value = 44.5945945945946
$ws.Range("F3").Value = 24
$ws.Range("G3").Value = 61.53846153846154

# Row 4: >10km
$ws.Range("A4").Value = ">10km"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 16.66666666666666
$ws.Range("D4").Value = 17
$ws.Range("E4").Value = 22.97297297297298
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = 10.25641025641026

# Row 5: clear label (no address group), new percentages/counts
$ws.Range("A5").Value = ""
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 11.11111111111111
$ws.Range("D5").Value = 7
$ws.Range("E5").Value = 9.45945945945946
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 5.128205128205128

# Remove old row 6 entirely (data now only spans through row 5)
$ws.Range("A6:G6").ClearContents()

# Update the used dimension range
$ws.Range("A1:G5").Select() | Out-Null
